$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC: set 27 cell values
$ws.Range("H19").Value = 900
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 900
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 900
$ws.Range("N19").Value = -1250
$ws.Range("H94").Value = 522525.12
$ws.Range("I94").Value = 602336.7
$ws.Range("J94").Value = 3750
$ws.Range("K94").Value = 602336.7
$ws.Range("L94").Value = 3750
$ws.Range("M94").Value = -601885.7
$ws.Range("N94").Value = -4652
$ws.Range("H137").Value = 3006.077
$ws.Range("I137").Value = 990.88
$ws.Range("J137").Value = 9723.4
$ws.Range("K137").Value = 2972.64
$ws.Range("L137").Value = 29170.2
$ws.Range("M137").Value = -422.6399999999999
$ws.Range("N137").Value = -34270.2
$ws.Range("H138").Value = 1949.32
$ws.Range("I138").Value = 802.7619
$ws.Range("J138").Value = 2779.5862
$ws.Range("K138").Value = 2408.2857
$ws.Range("L138").Value = 8338.758600000001
$ws.Range("M138").Value = 2731.7143
$ws.Range("N138").Value = -18618.7586
# ALC: clear 1 cells to blank
$ws.Range("M19").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
# ARM: set 18 cell values
$ws.Range("H88").Value = 12285626
$ws.Range("I88").Value = 33336342
$ws.Range("J88").Value = 2569910.8
$ws.Range("K88").Value = 33336342
$ws.Range("L88").Value = 2569910.8
$ws.Range("M88").Value = -33335936
$ws.Range("N88").Value = -2570722.8
$ws.Range("H91").Value = 12285626
$ws.Range("I91").Value = 33336342
$ws.Range("J91").Value = 2569910.8
$ws.Range("K91").Value = 33336342
$ws.Range("L91").Value = 2569910.8
$ws.Range("M91").Value = -33334938
$ws.Range("N91").Value = -2572718.8
$ws.Range("H106").Value = 47340.5
$ws.Range("J106").Value = 47340.5
$ws.Range("L106").Value = 47340.5
$ws.Range("N106").Value = -49864.5

$ws = $wb.Worksheets.Item("BSM")
# BSM: set 24 cell values
$ws.Range("H86").Value = 1936.1111
$ws.Range("I86").Value = 1890.0667
$ws.Range("J86").Value = 2166.3333
$ws.Range("K86").Value = 1890.0667
$ws.Range("L86").Value = 2166.3333
$ws.Range("M86").Value = -767.0667000000001
$ws.Range("N86").Value = -4412.3333
$ws.Range("H89").Value = 1936.1111
$ws.Range("I89").Value = 1890.0667
$ws.Range("J89").Value = 2166.3333
$ws.Range("K89").Value = 9450.333500000001
$ws.Range("L89").Value = 10831.6665
$ws.Range("M89").Value = -3834.333500000001
$ws.Range("N89").Value = -22063.6665
$ws.Range("H94").Value = 663.16
$ws.Range("I94").Value = 696.26666
$ws.Range("K94").Value = 696.26666
$ws.Range("M94").Value = -245.26666
$ws.Range("H99").Value = 4500
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 4500
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 4500
$ws.Range("N99").Value = -7496
# BSM: clear 1 cells to blank
$ws.Range("M99").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
# CRP: set 36 cell values
$ws.Range("H62").Value = 2853.4
$ws.Range("J62").Value = 3876.5
$ws.Range("L62").Value = 3876.5
$ws.Range("N62").Value = -5124.5
$ws.Range("H65").Value = 2853.4
$ws.Range("J65").Value = 3876.5
$ws.Range("L65").Value = 19382.5
$ws.Range("N65").Value = -25622.5
$ws.Range("H94").Value = 1622
$ws.Range("I94").Value = 1056
$ws.Range("J94").Value = 1716.3334
$ws.Range("K94").Value = 1056
$ws.Range("L94").Value = 1716.3334
$ws.Range("M94").Value = -605
$ws.Range("N94").Value = -2618.3334
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("H132").Value = 34996.023
$ws.Range("I132").Value = 1503.7778
$ws.Range("K132").Value = 4511.3334
$ws.Range("M132").Value = -1981.3334
$ws.Range("H134").Value = 218835.62
$ws.Range("I134").Value = 870.0645
$ws.Range("J134").Value = 669297.8
$ws.Range("K134").Value = 2610.1935
$ws.Range("L134").Value = 2007893.4
$ws.Range("M134").Value = -75.19349999999986
$ws.Range("N134").Value = -2012963.4
# CRP: clear 4 cells to blank
$ws.Range("M99").ClearContents()
$ws.Range("N99").ClearContents()
$ws.Range("M126").ClearContents()
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
# CUL: set 25 cell values
$ws.Range("H59").Value = 5000
$ws.Range("I59").Value = 0
$ws.Range("K59").Value = 0
$ws.Range("H87").Value = 34
$ws.Range("I87").Value = 34
$ws.Range("K87").Value = 102
$ws.Range("M87").Value = 1146
$ws.Range("H90").Value = 34
$ws.Range("I90").Value = 34
$ws.Range("K90").Value = 306
$ws.Range("M90").Value = 5934
$ws.Range("H113").Value = 28317300
$ws.Range("I113").Value = 50005572
$ws.Range("J113").Value = 2801687.5
$ws.Range("K113").Value = 150016716
$ws.Range("L113").Value = 8405062.5
$ws.Range("M113").Value = -150014546
$ws.Range("N113").Value = -8409402.5
$ws.Range("H131").Value = 7854.1875
$ws.Range("I131").Value = 11518.556
$ws.Range("J131").Value = 3142.8572
$ws.Range("K131").Value = 34555.66800000001
$ws.Range("L131").Value = 9428.571599999999
$ws.Range("M131").Value = -29515.66800000001
$ws.Range("N131").Value = -19508.5716
# CUL: clear 1 cells to blank
$ws.Range("M59").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
# GSM: set 6 cell values
$ws.Range("H102").Value = 1832.8182
$ws.Range("I102").Value = 1832.8182
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 1832.8182
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -210.8181999999999
# GSM: clear 1 cells to blank
$ws.Range("N102").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# LTW: set 73 cell values
$ws.Range("H7").Value = 2072.138
$ws.Range("I7").Value = 1841.125
$ws.Range("J7").Value = 3181
$ws.Range("K7").Value = 1841.125
$ws.Range("L7").Value = 3181
$ws.Range("M7").Value = -1729.125
$ws.Range("N7").Value = -3405
$ws.Range("H16").Value = 2296.7778
$ws.Range("I16").Value = 2241.6
$ws.Range("J16").Value = 2454.4285
$ws.Range("K16").Value = 2241.6
$ws.Range("L16").Value = 2454.4285
$ws.Range("M16").Value = -2071.6
$ws.Range("N16").Value = -2794.4285
$ws.Range("H51").Value = 32000
$ws.Range("J51").Value = 32000
$ws.Range("L51").Value = 32000
$ws.Range("N51").Value = -32956
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("H68").Value = 2458.077
$ws.Range("I68").Value = 2111
$ws.Range("J68").Value = 2675
$ws.Range("K68").Value = 2111
$ws.Range("L68").Value = 2675
$ws.Range("M68").Value = -1362
$ws.Range("N68").Value = -4173
$ws.Range("H71").Value = 2458.077
$ws.Range("I71").Value = 2111
$ws.Range("J71").Value = 2675
$ws.Range("K71").Value = 10555
$ws.Range("L71").Value = 13375
$ws.Range("M71").Value = -6811
$ws.Range("N71").Value = -20863
$ws.Range("H82").Value = 9261171
$ws.Range("I82").Value = 2028.7142
$ws.Range("J82").Value = 41668170
$ws.Range("K82").Value = 2028.7142
$ws.Range("L82").Value = 41668170
$ws.Range("M82").Value = -1667.7142
$ws.Range("N82").Value = -41668892
$ws.Range("H85").Value = 9261171
$ws.Range("I85").Value = 2028.7142
$ws.Range("J85").Value = 41668170
$ws.Range("K85").Value = 2028.7142
$ws.Range("L85").Value = 41668170
$ws.Range("M85").Value = -780.7141999999999
$ws.Range("N85").Value = -41670666
$ws.Range("H93").Value = 2085.5
$ws.Range("I93").Value = 1191.2727
$ws.Range("J93").Value = 2842.1538
$ws.Range("K93").Value = 1191.2727
$ws.Range("L93").Value = 2842.1538
$ws.Range("M93").Value = 56.72730000000001
$ws.Range("N93").Value = -5338.1538
$ws.Range("H100").Value = 2302.25
$ws.Range("I100").Value = 1920.6364
$ws.Range("J100").Value = 6500
$ws.Range("K100").Value = 1920.6364
$ws.Range("L100").Value = 6500
$ws.Range("M100").Value = -1379.6364
$ws.Range("N100").Value = -7582
$ws.Range("H126").Value = 2072.138
$ws.Range("I126").Value = 1841.125
$ws.Range("J126").Value = 3181
$ws.Range("K126").Value = 5523.375
$ws.Range("L126").Value = 9543
$ws.Range("M126").Value = -3053.375
$ws.Range("N126").Value = -14483
# LTW: clear 2 cells to blank
$ws.Range("N64").ClearContents()
$ws.Range("N67").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
# WVR: set 31 cell values
$ws.Range("H63").Value = 61329.332
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 61329.332
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 61329.332
$ws.Range("N63").Value = -62577.332
$ws.Range("H66").Value = 61329.332
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 61329.332
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 183987.996
$ws.Range("N66").Value = -190227.996
$ws.Range("H101").Value = 30118
$ws.Range("J101").Value = 30118
$ws.Range("L101").Value = 30118
$ws.Range("N101").Value = -36608
$ws.Range("H122").Value = 2381784.2
$ws.Range("I122").Value = 3572214
$ws.Range("J122").Value = 925
$ws.Range("K122").Value = 10716642
$ws.Range("L122").Value = 2775
$ws.Range("M122").Value = -10714192
$ws.Range("N122").Value = -7675
$ws.Range("H132").Value = 1688.6666
$ws.Range("I132").Value = 1536.9744
$ws.Range("K132").Value = 4610.9232
$ws.Range("M132").Value = -2080.9232
$ws.Range("H141").Value = 17644.889
$ws.Range("J141").Value = 17644.889
$ws.Range("L141").Value = 17644.889
$ws.Range("N141").Value = -28004.889
# WVR: clear 2 cells to blank
$ws.Range("M63").ClearContents()
$ws.Range("M66").ClearContents()
